# data formatting and cleaning update
# Remove the "(raw)" price rows for foods that also have a "(boiled)"/"(roasted)"
# preparation already listed, since the raw price duplicated/obsoleted those
# entries. Deleting the rows shifts everything below them up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row numbers (1-based, as currently laid out) that must be removed.
# Delete from the bottom up so earlier deletions don't renumber the
# rows we still need to remove.
$rowsToDelete = @(51, 50, 49, 32, 29, 20, 18, 17, 12, 9, 6, 2)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
